$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.318185925483704
$ws.Range("B1").Value = 2.540267944335938
$ws.Range("D1").Value = 1.591284871101379
$ws.Range("E1").Value = 0.9482159614562988
